# The deck's theme (theme1.xml, "Integral") is being swapped for the
# "Office Theme" palette that previously lived in theme2.xml (the notes
# master's theme). The PowerPoint object model doesn't expose a generic
# "swap these two theme parts" operation, but it does expose the live
# theme color scheme for editing via ThemeColorScheme.Colors(i).RGB
# (MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). Note RGB is a COLORREF-style integer (0xBBGGRR), so each
# target RRGGBB hex value below is byte-swapped before assignment.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000  # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444  # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED  # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244  # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70  # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305  # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95  # folHlink -> 954F72
